$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 7 (Ano 2025) with refreshed figures
$ws.Range("B7").Value = 2176806.51
$ws.Range("C7").Value = -51.81319521884107
$ws.Range("D7").Value = 2054
$ws.Range("E7").Value = 2054
$ws.Range("F7").Value = 1059.788953261928
$ws.Range("G7").Value = 9.37044006317571
